$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp in A1
$ws.Range("A1").Value = "Datos actualizados a 31 de Agosto de 2020 a las 18:43"

# Grecia's case counts grew past Croacia's, so the two swap places in the
# sorted (by total cases, descending) country list: Grecia now occupies the
# row that used to hold Croacia, and Croacia (unchanged data) shifts down
# into the row Grecia used to occupy.
$ws.Range("A91").Value = "Grecia"
$ws.Range("B91").Value = 10317
$ws.Range("C91").Value = 183
$ws.Range("D91").Value = 3804
$ws.Range("E91").Value = 6247
$ws.Range("F91").Value = 0
$ws.Range("G91").Value = 4
$ws.Range("H91").Value = 266

$ws.Range("A92").Value = "Croacia"
$ws.Range("B92").Value = 10269
$ws.Range("C92").Value = 146
$ws.Range("D92").Value = 7434
$ws.Range("E92").Value = 2649
$ws.Range("F92").Value = 0
$ws.Range("G92").Value = 2
$ws.Range("H92").Value = 186

# Estados Unidos
$ws.Range("B4").Value = 6185243
$ws.Range("C4").Value = 12007
$ws.Range("D4").Value = 3429638
$ws.Range("E4").Value = 2568229
$ws.Range("G4").Value = 152
$ws.Range("H4").Value = 187376

# India
$ws.Range("B6").Value = 3679411
$ws.Range("C6").Value = 60242
$ws.Range("D6").Value = 2832738
$ws.Range("E6").Value = 781269
$ws.Range("G6").Value = 787
$ws.Range("H6").Value = 65404

# Chile
$ws.Range("B13").Value = 411726
$ws.Range("C13").Value = 1752
$ws.Range("D13").Value = 383879
$ws.Range("E13").Value = 16558
$ws.Range("G13").Value = 45
$ws.Range("H13").Value = 11289

# Turquia
$ws.Range("B21").Value = 270133
$ws.Range("C21").Value = 1587
$ws.Range("D21").Value = 244926
$ws.Range("E21").Value = 18837
$ws.Range("G21").Value = 44
$ws.Range("H21").Value = 6370

# Alemania
$ws.Range("B23").Value = 244530
$ws.Range("C23").Value = 1235
$ws.Range("E23").Value = 17318
$ws.Range("G23").Value = 6
$ws.Range("H23").Value = 9370

# Emiratos Arabes Unidos
$ws.Range("B46").Value = 70231
$ws.Range("C46").Value = 541
$ws.Range("D46").Value = 60931
$ws.Range("E46").Value = 8916
$ws.Range("G46").Value = 2
$ws.Range("H46").Value = 384

# Suiza
$ws.Range("E61").Value = 4371
$ws.Range("G61").Value = 1
$ws.Range("H61").Value = 2006

# Chequia
$ws.Range("B74").Value = 24473
$ws.Range("C74").Value = 106
$ws.Range("D74").Value = 17772
$ws.Range("E74").Value = 6277
$ws.Range("G74").Value = 1
$ws.Range("H74").Value = 424

# Birmania
$ws.Range("B167").Value = 882
$ws.Range("C167").Value = 107
$ws.Range("D167").Value = 354
$ws.Range("E167").Value = 522
